$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")

# Set the shading type (type_shade, column L) for each building-use archetype row.
# Default value is "T1" for all rows, except SWIMMING and PARKING which get "T0".
$ws.Range("L2").Value = "T1"
$ws.Range("L3").Value = "T1"
$ws.Range("L4").Value = "T1"
$ws.Range("L5").Value = "T1"
$ws.Range("L6").Value = "T1"
$ws.Range("L7").Value = "T1"
$ws.Range("L8").Value = "T1"
$ws.Range("L9").Value = "T1"
$ws.Range("L10").Value = "T1"
$ws.Range("L11").Value = "T1"
$ws.Range("L12").Value = "T1"
$ws.Range("L13").Value = "T0"
$ws.Range("L14").Value = "T1"
$ws.Range("L15").Value = "T0"
$ws.Range("L16").Value = "T1"
$ws.Range("L17").Value = "T1"
$ws.Range("L18").Value = "T1"
$ws.Range("L19").Value = "T1"

# Match the active-cell selection left behind in the saved workbook.
$ws.Range("L1").Select()
